# Apply the edit: rotate columns G:J (paramName moves from G to J),
# unshare/refresh derived text-formula columns (X, AA, AB, AC, AF) and
# update AF's paramName reference from G to J, adjust column widths,
# and set the final selection to K1, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rotate columns G:J left for rows 1-15 ---
# new G = old H, new H = old I, new I = old J, new J = old G
$ws.Range("G1").Value = 'dataprintList'
$ws.Range("H1").Value = 'assumedDistrChoices'
$ws.Range("I1").Value = 'simXAxis_param'
$ws.Range("J1").Value = 'paramName'
$ws.Range("G2").Value = 'intPrintHelper'
$ws.Range("H2").Value = 'list("Bernoulli-Pi","Bernoulli-Logit", "Bernoulli-Logit-X")'
$ws.Range("I2").Value = '$ \\tilde{E}(y) = \\tilde{\\pi} = \\tilde{Pr}(Y=1)$'
$ws.Range("J2").Value = 'Pi'
$ws.Range("G3").Value = 'intPrintHelper'
$ws.Range("H3").Value = 'list("Bernoulli-Logit","Bernoulli-Pi", "Bernoulli-Logit-X")'
$ws.Range("I3").Value = '$ \\tilde{E}(y) =\\tilde{\\pi} = \\tilde{Pr}(Y=1)$'
$ws.Range("J3").Value = 'Beta'
$ws.Range("G4").Value = 'intPrintHelper'
$ws.Range("H4").Value = 'list("Bernoulli-Logit-X","Bernoulli-Pi","Bernoulli-Logit")'
$ws.Range("I4").Value = '$ \\tilde{E}(y) =\\tilde{\\pi} = \\tilde{Pr}(Y=1)$'
$ws.Range("J4").Value = 'Beta'
$ws.Range("G5").Value = 'decPrintHelper'
$ws.Range("H5").Value = 'list("Stylized-Normal","Stylized-Normal-X")'
$ws.Range("I5").Value = '$ \\tilde{E}(y) =\\tilde{\\mu} = \\bar{Y}$'
$ws.Range("J5").Value = 'Beta'
$ws.Range("G6").Value = 'decPrintHelper'
$ws.Range("H6").Value = 'list("Stylized-Normal-X","Stylized-Normal")'
$ws.Range("I6").Value = '$ \\tilde{E}(y) =\\tilde{\\mu} = \\bar{Y}$'
$ws.Range("J6").Value = 'Beta'
$ws.Range("G7").Value = 'decPrintHelper'
$ws.Range("H7").Value = 'list("Normal-X", "Stylized-Normal-X","Stylized-Normal")'
$ws.Range("I7").Value = '$ \\tilde{E}(y) =\\tilde{\\mu} = \\bar{Y}$'
$ws.Range("J7").Value = 'Beta/Sigma'
$ws.Range("G8").Value = 'decPrintHelper'
$ws.Range("H8").Value = 'list("Log-Normal","Stylized-Normal","Stylized-Normal-X")'
$ws.Range("I8").Value = '$ \\tilde{E}(y)$'
$ws.Range("J8").Value = 'Beta'
$ws.Range("G9").Value = 'decPrintHelper'
$ws.Range("H9").Value = 'list("Log-Normal-X","Stylized-Normal","Stylized-Normal-X")'
$ws.Range("I9").Value = '$ \\tilde{E}(y)$'
$ws.Range("J9").Value = 'Beta'
$ws.Range("G10").Value = 'intPrintHelper'
$ws.Range("H10").Value = 'list("Poisson", "Poisson-Exp", "Poisson-Exp-X", "Stylized-Normal", "Stylized-Normal-X", "Normal-X")'
$ws.Range("I10").Value = '$ \\tilde{E}(y) =\\tilde{\\lambda} = \\bar{Y}$'
$ws.Range("J10").Value = 'Lambda'
$ws.Range("G11").Value = 'intPrintHelper'
$ws.Range("H11").Value = 'list("Poisson", "Poisson-Exp", "Poisson-Exp-X", "Stylized-Normal", "Stylized-Normal-X", "Normal-X")'
$ws.Range("I11").Value = '$ \\tilde{E}(y) =\\tilde{\\lambda} = \\bar{Y}$'
$ws.Range("J11").Value = 'Beta'
$ws.Range("G12").Value = 'intPrintHelper'
$ws.Range("H12").Value = 'list("Poisson", "Poisson-Exp", "Poisson-Exp-X", "Stylized-Normal", "Stylized-Normal-X", "Normal-X")'
$ws.Range("I12").Value = '$ \\tilde{E}(y) =\\tilde{\\lambda} = \\bar{Y}$'
$ws.Range("J12").Value = 'Beta'
$ws.Range("G13").Value = 'decPrintHelper'
$ws.Range("H13").Value = 'list("Exponential", "Stylized-Normal","Stylized-Normal-X","Log-Normal", "Log-Normal-X","Exponential-Exp", "Exponential-Exp-X")'
$ws.Range("I13").Value = '$ \\tilde{E}(y)$'
$ws.Range("J13").Value = 'Lambda'
$ws.Range("G14").Value = 'decPrintHelper'
$ws.Range("H14").Value = 'list("Exponential-Exp","Stylized-Normal","Stylized-Normal-X","Log-Normal", "Log-Normal-X","Exponential",  "Exponential-Exp-X")'
$ws.Range("I14").Value = '$ \\tilde{E}(y)$'
$ws.Range("J14").Value = 'Beta'
$ws.Range("G15").Value = 'decPrintHelper'
$ws.Range("H15").Value = 'list( "Exponential-Exp-X","Stylized-Normal","Stylized-Normal-X","Log-Normal", "Log-Normal-X","Exponential", "Exponential-Exp")'
$ws.Range("I15").Value = '$ \\tilde{E}(y)$'
$ws.Range("J15").Value = 'Beta'

# --- 2) Column width adjustments (G/H/I/J) ---
$ws.Columns.Item(7).ColumnWidth = 14.5703125
$ws.Columns.Item(8).ColumnWidth = 14.5703125
$ws.Columns.Item(9).ColumnWidth = 8.43
$ws.Columns.Item(10).ColumnWidth = 15.5703125

# --- 3) Refresh derived formula columns X, AA, AB, AC (unchanged logic, re-enter
#        so the workbook reflects them as individually edited, not shared) ---
$ws.Range("X2").Formula = '=IF(E2=1,"c()", IF(E2=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA2").Formula = '=B2&"PlotDistr"'
$ws.Range("AB2").Formula = '=B2&"Draws"'
$ws.Range("AC2").Formula = '=B2&"Latex"'
$ws.Range("X3").Formula = '=IF(E3=1,"c()", IF(E3=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA3").Formula = '=B3&"PlotDistr"'
$ws.Range("AB3").Formula = '=B3&"Draws"'
$ws.Range("AC3").Formula = '=B3&"Latex"'
$ws.Range("X4").Formula = '=IF(E4=1,"c()", IF(E4=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA4").Formula = '=B4&"PlotDistr"'
$ws.Range("AB4").Formula = '=B4&"Draws"'
$ws.Range("AC4").Formula = '=B4&"Latex"'
$ws.Range("X5").Formula = '=IF(E5=1,"c()", IF(E5=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA5").Formula = '=B5&"PlotDistr"'
$ws.Range("AB5").Formula = '=B5&"Draws"'
$ws.Range("AC5").Formula = '=B5&"Latex"'
$ws.Range("X6").Formula = '=IF(E6=1,"c()", IF(E6=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA6").Formula = '=B6&"PlotDistr"'
$ws.Range("AB6").Formula = '=B6&"Draws"'
$ws.Range("AC6").Formula = '=B6&"Latex"'
$ws.Range("X7").Formula = '=IF(E7=1,"c()", IF(E7=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA7").Formula = '=B7&"PlotDistr"'
$ws.Range("AB7").Formula = '=B7&"Draws"'
$ws.Range("AC7").Formula = '=B7&"Latex"'
$ws.Range("X8").Formula = '=IF(E8=1,"c()", IF(E8=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA8").Formula = '=B8&"PlotDistr"'
$ws.Range("AB8").Formula = '=B8&"Draws"'
$ws.Range("AC8").Formula = '=B8&"Latex"'
$ws.Range("X9").Formula = '=IF(E9=1,"c()", IF(E9=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA9").Formula = '=B9&"PlotDistr"'
$ws.Range("AB9").Formula = '=B9&"Draws"'
$ws.Range("AC9").Formula = '=B9&"Latex"'
$ws.Range("X10").Formula = '=IF(E10=1,"c()", IF(E10=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA10").Formula = '=B10&"PlotDistr"'
$ws.Range("AB10").Formula = '=B10&"Draws"'
$ws.Range("AC10").Formula = '=B10&"Latex"'
$ws.Range("X11").Formula = '=IF(E11=1,"c()", IF(E11=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA11").Formula = '=B11&"PlotDistr"'
$ws.Range("AB11").Formula = '=B11&"Draws"'
$ws.Range("AC11").Formula = '=B11&"Latex"'
$ws.Range("X12").Formula = '=IF(E12=1,"c()", IF(E12=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA12").Formula = '=B12&"PlotDistr"'
$ws.Range("AB12").Formula = '=B12&"Draws"'
$ws.Range("AC12").Formula = '=B12&"Latex"'
$ws.Range("X13").Formula = '=IF(E13=1,"c()", IF(E13=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA13").Formula = '=B13&"PlotDistr"'
$ws.Range("AB13").Formula = '=B13&"Draws"'
$ws.Range("AC13").Formula = '=B13&"Latex"'
$ws.Range("X14").Formula = '=IF(E14=1,"c()", IF(E14=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA14").Formula = '=B14&"PlotDistr"'
$ws.Range("AB14").Formula = '=B14&"Draws"'
$ws.Range("AC14").Formula = '=B14&"Latex"'
$ws.Range("X15").Formula = '=IF(E15=1,"c()", IF(E15=3,"c(""Beta0"", ""Beta1"", ""Beta2"")","c(""Beta0"", ""Beta1"", ""Beta2"",""Sigma"")"))'
$ws.Range("AA15").Formula = '=B15&"PlotDistr"'
$ws.Range("AB15").Formula = '=B15&"Draws"'
$ws.Range("AC15").Formula = '=B15&"Latex"'

# --- 4) AF column: paramName now read from J (was G) ---
$ws.Range("AF2").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD2&", likelihoodFun = "&AE2&" , paramName = """&J2&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N2&""")}"'
$ws.Range("AF3").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD3&", likelihoodFun = "&AE3&" , paramName = """&J3&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N3&""")}"'
$ws.Range("AF4").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD4&", likelihoodFun = "&AE4&" , paramName = """&J4&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N4&""")}"'
$ws.Range("AF5").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD5&", likelihoodFun = "&AE5&" , paramName = """&J5&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N5&""")}"'
$ws.Range("AF6").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD6&", likelihoodFun = "&AE6&" , paramName = """&J6&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N6&""")}"'
$ws.Range("AF7").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD7&", likelihoodFun = "&AE7&" , paramName = """&J7&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N7&""")}"'
$ws.Range("AF8").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD8&", likelihoodFun = "&AE8&" , paramName = """&J8&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N8&""")}"'
$ws.Range("AF9").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD9&", likelihoodFun = "&AE9&" , paramName = """&J9&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N9&""")}"'
$ws.Range("AF10").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD10&", likelihoodFun = "&AE10&" , paramName = """&J10&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N10&""")}"'
$ws.Range("AF11").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD11&", likelihoodFun = "&AE11&" , paramName = """&J11&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N11&""")}"'
$ws.Range("AF12").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD12&", likelihoodFun = "&AE12&" , paramName = """&J12&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N12&""")}"'
$ws.Range("AF13").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD13&", likelihoodFun = "&AE13&" , paramName = """&J13&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N13&""")}"'
$ws.Range("AF14").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD14&", likelihoodFun = "&AE14&" , paramName = """&J14&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N14&""")}"'
$ws.Range("AF15").Formula = '="function(outcome, xVals, margNum){MLEstimator(outcome = outcome, chartDomain = "&AD15&", likelihoodFun = "&AE15&" , paramName = """&J15&""",  xVals = xVals, margNum = margNum, "&"optimMethod = """&N15&""")}"'

# --- 5) Final selection ---
$ws.Range("K1").Select() | Out-Null
